# Update the "想去人数" (interest count, column F) figures for a handful of
# events in the "展览" and "全部类型" sheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 64
$ws1.Range("F4").Value  = 3673
$ws1.Range("F5").Value  = 2247
$ws1.Range("F8").Value  = 4
$ws1.Range("F10").Value = 96
$ws1.Range("F12").Value = 1358
$ws1.Range("F14").Value = 2063
$ws1.Range("F15").Value = 148

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 64
$ws4.Range("F4").Value  = 3673
$ws4.Range("F5").Value  = 2247
$ws4.Range("F8").Value  = 4
$ws4.Range("F11").Value = 96
$ws4.Range("F15").Value = 1358
$ws4.Range("F17").Value = 2063
$ws4.Range("F18").Value = 148
